$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet after Sheet1 and name it Sheet2
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Populate the new sheet
$ws2.Range("A1").Value = "columnname"
$ws2.Range("A2").Value = "otherimportantdata"

# Sheet2 becomes the active/selected tab, with A3 selected (matching Sheet1's prior selection)
$ws2.Activate() | Out-Null
$ws2.Range("A3").Select() | Out-Null
